# Lecture partielle de l'EDT M1 MIAGE.
# Shift every scheduled date in column A forward by 1096 days (same
# month/day, 3 years later) and update the French weekday label in the
# adjacent column B so it still matches the (now different) day of week
# for that calendar date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,4,8,12,15,18,21,23,26,28,31,35,38,42,48,51,54,58,60,63)
$dayNames = @{
    2  = "mercredi"
    4  = "jeudi"
    8  = "jeudi"
    12 = "jeudi"
    15 = "jeudi"
    18 = "jeudi"
    21 = "lundi"
    23 = "jeudi"
    26 = "lundi"
    28 = "jeudi"
    31 = "lundi"
    35 = "mardi"
    38 = "jeudi"
    42 = "vendredi"
    48 = "mardi"
    51 = "mercredi"
    54 = "vendredi"
    58 = "mercredi"
    60 = "vendredi"
    63 = "mercredi"
}

foreach ($r in $rows) {
    $aCell = $ws.Range("A$r")
    $aCell.Value = $aCell.Value2 + 1096.0
    $ws.Range("B$r").Value = $dayNames[$r]
}
